$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 13:52"

# Alemania (row 8) - refreshed totals
$ws.Range("B8").Value = 153393
$ws.Range("C8").Value = 264
$ws.Range("E8").Value = 41018

# Iran (row 11) - refreshed totals
$ws.Range("B11").Value = 88194
$ws.Range("C11").Value = 1168
$ws.Range("D11").Value = 66599
$ws.Range("E11").Value = 16021
$ws.Range("F11").Value = 3121
$ws.Range("G11").Value = 93
$ws.Range("H11").Value = 5574

# Republica de Yibuti (row 87) - refreshed totals
$ws.Range("B87").Value = 999
$ws.Range("C87").Value = 13
$ws.Range("D87").Value = 330
$ws.Range("E87").Value = 667

# Brunei (row 136) - refreshed totals
$ws.Range("D136").Value = 120
$ws.Range("E136").Value = 17

# Liberia jumps up in the ranking to sit right after Etiopia (row 141),
# pushing Maldivas / Trinidad yTobago / Guayana Francesa down one row each.
# Row 145 (Aruba) is unaffected.
$ws.Range("A141").Value = "Liberia"
$ws.Range("B141").Value = 117
$ws.Range("C141").Value = 16
$ws.Range("D141").Value = 25
$ws.Range("E141").Value = 84
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 8

$ws.Range("A142").Value = "Maldivas"
$ws.Range("B142").Value = 116
$ws.Range("C142").Value = 8
$ws.Range("D142").Value = 16
$ws.Range("E142").Value = 100
$ws.Range("F142").Value = 2
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

$ws.Range("A143").Value = "Trinidad yTobago"
$ws.Range("B143").Value = 115
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 48
$ws.Range("E143").Value = 59
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 8

$ws.Range("A144").Value = "Guayana Francesa"
$ws.Range("B144").Value = 107
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 84
$ws.Range("E144").Value = 22
$ws.Range("F144").Value = 1
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 1

# Suazilandia jumps up in the ranking to sit right after Mongolia (row 174),
# pushing Malaui / Guam down one row each. Row 177 (Zimbabue) is unaffected.
$ws.Range("A174").Value = "Suazilandia"
$ws.Range("B174").Value = 36
$ws.Range("C174").Value = 5
$ws.Range("D174").Value = 10
$ws.Range("E174").Value = 25
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 1

$ws.Range("A175").Value = "Malaui"
$ws.Range("B175").Value = 33
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 3
$ws.Range("E175").Value = 27
$ws.Range("F175").Value = 1
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 3

$ws.Range("A176").Value = "Guam"
$ws.Range("B176").Value = 32
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 31
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 1
